# Updates the cryptos list worksheet with refreshed price/volume data
# pulled on Fri Jul 19 12:50:04 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text can never be mistaken for a number (percentages,
# names, links, multi-dot "thousand-grouped" prices) -- plain assignment
# keeps their existing (default) cell style untouched.
$textUpdates = @{
    'D2' = '64.357.41'
    'E2' = '  -0.79%  '
    'D3' = '3.413.95'
    'E3' = '  -1.56%  '
    'E5' = '  -0.24%  '
    'E6' = '  +1.79%  '
    'E7' = '  +0.07%  '
    'D8' = '3.412.93'
    'E8' = '  -1.54%  '
    'E9' = '  -4.45%  '
    'E10' = '  +1.15%  '
    'E11' = '  -1.98%  '
    'E12' = '  -4.44%  '
    'D13' = '4.001.38'
    'E13' = '  -1.47%  '
    'E14' = '  +0.06%  '
    'E15' = '  -2.66%  '
    'E16' = '  +0.95%  '
    'D17' = '64.353.29'
    'E17' = '  -0.89%  '
    'D18' = '3.413.25'
    'E18' = '  -1.61%  '
    'E19' = '  -1.10%  '
    'E20' = '  -2.51%  '
    'E21' = '  -1.54%  '
    'E22' = '  -2.22%  '
    'E23' = '  +0.12%  '
    'E24' = '  -3.07%  '
    'E25' = '  -3.13%  '
    'E26' = '  -3.91%  '
    'E27' = '  -4.65%  '
    'E28' = '  -0.79%  '
    'E29' = '  -0.18%  '
    'E30' = '  -1.80%  '
    'E31' = '  -3.95%  '
    'E32' = '  -0.83%  '
    'E33' = '  +0.01%  '
    'E34' = '  -1.89%  '
    'E35' = '  -0.46%  '
    'E36' = '  -7.35%  '
    'E37' = '  -1.28%  '
    'E38' = '  +6.33%  '
    'E39' = '  -3.40%  '
    'B40' = 'EnergySwap'
    'C40' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E40' = '  -3.25%  '
    'B41' = 'Hedera'
    'C41' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E41' = '  -4.50%  '
    'D42' = '2.770.59'
    'E42' = '  -4.13%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'E43' = '  -1.46%  '
    'B44' = 'OKB'
    'C44' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'E44' = '  -1.01%  '
    'E45' = '  -0.76%  '
    'E46' = '  -3.75%  '
    'E47' = '  -2.24%  '
    'E48' = '  -1.46%  '
    'E49' = '  +2.33%  '
    'E50' = '  -3.23%  '
    'E51' = '  -1.71%  '
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price cells that look like plain decimals (e.g. "574.42") would be
# auto-coerced to numbers by a bare .Value assignment, which would both
# change the stored type and drop significant trailing zeros ("1.00" ->
# 1). Force them to text via a temporary "@" number format, then restore
# the cell's original (default/"Normal") style so formatting is unchanged.
$numericLookingUpdates = @{
    'D5' = '574.42'
    'D6' = '163.95'
    'D9' = '0.556'
    'D11' = '0.120'
    'D12' = '0.421'
    'D15' = '26.94'
    'D16' = '0.0000173'
    'D19' = '6.16'
    'D20' = '13.50'
    'D21' = '375.64'
    'D22' = '7.82'
    'D23' = '1.00'
    'D24' = '70.38'
    'D25' = '0.515'
    'D27' = '9.47'
    'D29' = '0.999'
    'D30' = '6.09'
    'D31' = '1.40'
    'D33' = '0.999'
    'D34' = '22.87'
    'D35' = '7.03'
    'D36' = '1.48'
    'D37' = '159.07'
    'D38' = '0.861'
    'D39' = '1.83'
    'D40' = '25.91'
    'D41' = '0.0722'
    'D43' = '6.49'
    'D44' = '42.63'
    'D45' = '25.79'
    'D46' = '4.37'
    'D47' = '0.0304'
    'D48' = '2.43'
    'D49' = '329.94'
}

foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    $cell.Style = "Normal"
}
